$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 50, pushing the existing blank/summary rows down by one.
$ws.Rows("50:50").Insert()

# Fill the new row 50 with a new data entry (2014-03-06, 08:00-12:00 -> 240 min -> 4h)
$ws.Cells.Item(50, 1).Value = 2014
$ws.Cells.Item(50, 2).Value = 3
$ws.Cells.Item(50, 3).Value = 6
$ws.Cells.Item(50, 4).Value = 0.33333333333333331
$ws.Cells.Item(50, 5).Value = 0.5

$ws.Cells.Item(50, 4).NumberFormat = $ws.Cells.Item(49, 4).NumberFormat
$ws.Cells.Item(50, 5).NumberFormat = $ws.Cells.Item(49, 5).NumberFormat
$ws.Cells.Item(50, 6).NumberFormat = $ws.Cells.Item(49, 6).NumberFormat
$ws.Cells.Item(50, 7).NumberFormat = $ws.Cells.Item(49, 7).NumberFormat

$ws.Cells.Item(50, 6).Formula = "=(E50-D50)*24*60"
$ws.Cells.Item(50, 7).Formula = "=F50/60"

# Fix the "sum [min]" formula so it sums through the new blank row 51 as before.
$ws.Cells.Item(52, 6).Formula = "=SUM(F2:F51)"

# Restore the selection/active cell as recorded after the edit.
$ws.Range("A51").Select()
